$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 4
}

$newB = @{
    2 = 546
    3 = 502
    4 = 477
    5 = 459
    6 = 408
    7 = 382
    8 = 373
    9 = 374
    10 = 355
    11 = 344
    12 = 336
    13 = 307
    14 = 292
    15 = 278
    16 = 267
    17 = 252
    18 = 232
    19 = 208
    20 = 169
    21 = 152
    22 = 148
    23 = 134
    24 = 127
    25 = 134
    26 = 142
    27 = 136
    28 = 131
    29 = 133
    30 = 129
    32 = 114
    33 = 137
    34 = 130
    35 = 129
    36 = 123
    37 = 118
    38 = 131
    39 = 148
}

foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $newB[$r]
}
